$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix duplicated "类" typo in the column headers (row 1)
$ws.Range("B1").Value = "其他服务类居民消费价格指数(上年=100)"
$ws.Range("D1").Value = "其他用品类居民消费价格指数(上年=100)"

# Seed the new year rows (7 & 8) with the same look as the existing
# year-label cells (bold, bordered, centered) by copying row 6's formatting.
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A8").PasteSpecial(-4122)

# Row 7: 2021年
$ws.Range("A7").Value = "2021年"
$ws.Range("B7").Value = 97.90000000000001
$ws.Range("C7").Value = 98.7
$ws.Range("D7").Value = 99.5

# Row 8: 2022年 (only the C column value is known/published)
$ws.Range("A8").Value = "2022年"
$ws.Range("C8").Value = 101.6
